$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tags")

$data = @(
    @(18350, "cell-type", 0),
    @(18650, "cell-type", 1),
    @(20700, "cell-type", 2),
    @(21700, "cell-type", 3),
    @(26650, "cell-type", 4)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Range("C20").Select()
